# Add a new tea entry ("Gyokuro") as row 20 of the Tea worksheet.
#
# The existing data block is rows 2-19; row 19 ("Kuki Cha" / Russian
# Samovar / Mate) is the last data row and carries the per-column number
# formatting (s="1"/"2"/"3") that the new row should inherit. We copy row
# 19's formatting down into row 20 first, then overwrite the new row's
# values/formulas - this keeps the shared "Theanine Share" (D) and
# "Average Elevation" (H) formulas extending cleanly into the new row,
# matching the same style indices used throughout the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone formatting (styles/number formats) of the last existing data row
# down into the new row.
$ws.Range("A19:J19").Copy()
$ws.Range("A20:J20").PasteSpecial(-4122)

# Gyokuro - Japanese shade-grown green tea.
$ws.Range("A20").Value = "Gyokuro"
$ws.Range("B20").Value = 30.0
$ws.Range("C20").Value = 40.0
$ws.Range("D20").Formula = "=B20/(C20+B20)"
$ws.Range("E20").Value = "Japan"
$ws.Range("F20").Value = 400.0
$ws.Range("G20").Value = 600.0
$ws.Range("H20").Formula = "=(F20+G20)/2"
$ws.Range("I20").Value = "green"
$ws.Range("J20").Value = "green"
